$wb = $excel.ActiveWorkbook

# --- 1. Insert the new "2022-Q1" sheet right before "总计" ---
$total = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($total)
$newSheet.Name = "2022-Q1"

# $total was a positional reference; after the insert it now points at the
# newly-added sheet, so re-resolve it by name to get back to "总计".
$total = $wb.Worksheets.Item("总计")

# Copy header-row look & feel (bold/border/center-top) from the "2021-Q4" sheet,
# which has the identical 7-column layout already.
$src = $wb.Worksheets.Item("2021-Q4")
$src.Range("A1:H2").Copy()
$newSheet.Range("A1").PasteSpecial(-4122)

# --- header row (B1:H1) ---
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# --- data row (A2:H2) ---
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'159962"
$newSheet.Range("B2").Style = "Normal"
$newSheet.Range("C2").Value = "华夏中证四川国企改革ETF"
$newSheet.Range("D2").Value = "'0.49"
$newSheet.Range("D2").Style = "Normal"
$newSheet.Range("E2").Value = "'95.82"
$newSheet.Range("E2").Style = "Normal"
$newSheet.Range("F2").Value = "'3.00"
$newSheet.Range("F2").Style = "Normal"
$newSheet.Range("G2").Value = "'0.0147"
$newSheet.Range("G2").Style = "Normal"
$newSheet.Range("H2").Value = 9

# --- 2. Update "总计": insert a fresh row 2 for 2022-Q1, push the rest down ---
$total.Rows.Item(2).Insert()
$total.Range("B2:D2").Style = "Normal"

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.01

# Renumber the index column for the rows that shifted down
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3

Write-Host "done"
